# Update "International Ever Green_2025-11-24.xlsx"
# - Append new order-line rows (42-71) to the "Orders" sheet, mirroring the
#   existing layout (PackageID in A only for the first line of a new
#   package, FlowerName in C, Number-as-text in F). All of A/C/F in this
#   sheet are text-typed (even the numeric-looking ones), matching the
#   pre-existing rows 2-41, so numeric-looking entries are force-formatted
#   as Text before assignment to avoid Excel's automatic number coercion.
# - Extend the ignoredErrors ("number stored as text") hint to the new range.
# - Append the new per-package digit codes to the "Summary" sheet G2 cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Orders")
$ws2 = $wb.Worksheets.Item("Summary")

$ws1.Cells.Item(42,3).Value = "816_山里红_undefined_undefined_1bunch"
$ws1.Cells.Item(42,6).NumberFormat = "@"
$ws1.Cells.Item(42,6).Value = "5"
$ws1.Cells.Item(43,1).NumberFormat = "@"
$ws1.Cells.Item(43,1).Value = "7"
$ws1.Cells.Item(43,3).Value = "653_大丽花 黑_undefined_undefined_5stems"
$ws1.Cells.Item(43,6).NumberFormat = "@"
$ws1.Cells.Item(43,6).Value = "8"
$ws1.Cells.Item(44,3).Value = "669_大丽花 红_undefined_undefined_5stems"
$ws1.Cells.Item(44,6).NumberFormat = "@"
$ws1.Cells.Item(44,6).Value = "10"
$ws1.Cells.Item(45,3).Value = "651_大丽花 奶油桃子_undefined_undefined_5stems"
$ws1.Cells.Item(45,6).NumberFormat = "@"
$ws1.Cells.Item(45,6).Value = "15"
$ws1.Cells.Item(46,3).Value = "314_松虫草花边黑_scabiosa_undefined_1bunch"
$ws1.Cells.Item(46,6).NumberFormat = "@"
$ws1.Cells.Item(46,6).Value = "7"
$ws1.Cells.Item(47,1).NumberFormat = "@"
$ws1.Cells.Item(47,1).Value = "8"
$ws1.Cells.Item(47,3).Value = "647_海棠果红_undefined_undefined_1bunch"
$ws1.Cells.Item(47,6).NumberFormat = "@"
$ws1.Cells.Item(47,6).Value = "6"
$ws1.Cells.Item(48,3).Value = "647_海棠果红_undefined_undefined_1bunch"
$ws1.Cells.Item(48,6).NumberFormat = "@"
$ws1.Cells.Item(48,6).Value = "5"
$ws1.Cells.Item(49,3).Value = "418_松虫草白_scabiosa white_undefined_1bunch"
$ws1.Cells.Item(49,6).NumberFormat = "@"
$ws1.Cells.Item(49,6).Value = "10"
$ws1.Cells.Item(50,3).Value = "512_松虫草粉_scabiosa pink_undefined_1bunch"
$ws1.Cells.Item(50,6).NumberFormat = "@"
$ws1.Cells.Item(50,6).Value = "6"
$ws1.Cells.Item(51,3).Value = "419_松虫草红_scabiosa watermelon_undefined_1bunch"
$ws1.Cells.Item(51,6).NumberFormat = "@"
$ws1.Cells.Item(51,6).Value = "12"
$ws1.Cells.Item(52,3).Value = "480_蝴蝶洋牡丹红_butterfly  Ranunculus_undefined_1bunch"
$ws1.Cells.Item(52,6).NumberFormat = "@"
$ws1.Cells.Item(52,6).Value = "6"
$ws1.Cells.Item(53,3).Value = "586_洋牡丹白_undefined_undefined_1bunch"
$ws1.Cells.Item(53,6).NumberFormat = "@"
$ws1.Cells.Item(53,6).Value = "10"
$ws1.Cells.Item(54,3).Value = "590_洋牡丹粉_undefined_undefined_1bunch"
$ws1.Cells.Item(54,6).NumberFormat = "@"
$ws1.Cells.Item(54,6).Value = "5"
$ws1.Cells.Item(55,3).Value = "585_洋牡丹红_undefined_undefined_1bunch"
$ws1.Cells.Item(55,6).NumberFormat = "@"
$ws1.Cells.Item(55,6).Value = "5"
$ws1.Cells.Item(56,1).NumberFormat = "@"
$ws1.Cells.Item(56,1).Value = "9"
$ws1.Cells.Item(56,3).Value = "462_五针松_undefined_undefined_1bunch"
$ws1.Cells.Item(56,6).NumberFormat = "@"
$ws1.Cells.Item(56,6).Value = "10"
$ws1.Cells.Item(57,3).Value = "359_蓝梦叶_undefined_undefined_1bunch"
$ws1.Cells.Item(57,6).NumberFormat = "@"
$ws1.Cells.Item(57,6).Value = "5"
$ws1.Cells.Item(58,3).Value = "411_紫罗兰白_violet white_undefined_1bunch"
$ws1.Cells.Item(58,6).NumberFormat = "@"
$ws1.Cells.Item(58,6).Value = "15"
$ws1.Cells.Item(59,3).Value = "412_紫罗兰粉_violet pink_undefined_1bunch"
$ws1.Cells.Item(59,6).NumberFormat = "@"
$ws1.Cells.Item(59,6).Value = "10"
$ws1.Cells.Item(60,3).Value = "506_紫罗兰香槟色_violet champagne_undefined_1bunch"
$ws1.Cells.Item(60,6).NumberFormat = "@"
$ws1.Cells.Item(60,6).Value = "5"
$ws1.Cells.Item(61,3).Value = "508_风铃花白色_Canterbury Bells " + [char]10 + "white_undefined_1bunch"
$ws1.Cells.Item(61,6).NumberFormat = "@"
$ws1.Cells.Item(61,6).Value = "5"
$ws1.Cells.Item(62,3).Value = "396_米花 白_rice flower white_undefined_1bunch"
$ws1.Cells.Item(62,6).NumberFormat = "@"
$ws1.Cells.Item(62,6).Value = "5"
$ws1.Cells.Item(63,3).Value = "387_洋甘菊_Chamomile_undefined_1bunch"
$ws1.Cells.Item(63,6).NumberFormat = "@"
$ws1.Cells.Item(63,6).Value = "5"
$ws1.Cells.Item(64,3).Value = "769_菟葵绿铃铛_undefined_undefined_undefinedundefined"
$ws1.Cells.Item(64,6).NumberFormat = "@"
$ws1.Cells.Item(64,6).Value = "10"
$ws1.Cells.Item(65,3).Value = "753_蝴蝶洋牡丹黄_butterfly  Ranunculus_undefined_1bunch"
$ws1.Cells.Item(65,6).NumberFormat = "@"
$ws1.Cells.Item(65,6).Value = "15"
$ws1.Cells.Item(66,3).Value = "773_格桑花白_undefined_undefined_1bunch"
$ws1.Cells.Item(66,6).NumberFormat = "@"
$ws1.Cells.Item(66,6).Value = "8"
$ws1.Cells.Item(67,3).Value = "405_小飞燕浅蓝_ delphinium ballkleid" + [char]10 + "dark blue_undefined_1bunch"
$ws1.Cells.Item(67,6).NumberFormat = "@"
$ws1.Cells.Item(67,6).Value = "10"
$ws1.Cells.Item(68,3).Value = "468_水仙百合_Alstroemeria_undefined_1bunch"
$ws1.Cells.Item(68,6).NumberFormat = "@"
$ws1.Cells.Item(68,6).Value = "5"
$ws1.Cells.Item(69,3).Value = "578_腊梅粉_wax pink_undefined_1bunch"
$ws1.Cells.Item(69,6).NumberFormat = "@"
$ws1.Cells.Item(69,6).Value = "15"
$ws1.Cells.Item(70,3).Value = "794_小菊罗西香槟_undefined_undefined_1bunch"
$ws1.Cells.Item(70,6).NumberFormat = "@"
$ws1.Cells.Item(70,6).Value = "10"
$ws1.Cells.Item(71,3).Value = "794_小菊罗西香槟_undefined_undefined_1bunch"


# Keep the "number stored as text" ignored-error hint in sync with the
# newly-populated range (A1:L71), matching the existing convention used
# for A1:L41.
$ws1.Range("A1:L71").Errors.Item(1).Ignore = $true

# Summary sheet: G2 accumulates one digit-group per order line across the
# whole sheet; append the digit groups contributed by the new rows above.
$ws2.Range("G2").NumberFormat = "@"
$ws2.Cells.Item(2,7).Value = "05515555251065256655552532151582255551055555125515655810157651061261055105151055551015810515100"
